# Mark "89ec4140-7020-4012-9fe1-624c2b8a2ebb" as ready for handoff, with a fresh
# handoff timestamp, across the Overview/zh-cn/de-de sheets.

$wb = $excel.ActiveWorkbook

$status = "Ready for handoff"
$newHandoffDateTime = "2016-03-21 18:47:36"
$newHandoffDateTimeZhCn = "2016-03-21 18:47:32"

# Overview sheet: row 3 is the 89ec4140... file
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status
$wsOverview.Range("D3").Value = $newHandoffDateTime

# zh-cn sheet: row 3 is the 89ec4140... file
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("E3").Value = $newHandoffDateTimeZhCn

# de-de sheet: row 3 is the 89ec4140... file
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("E3").Value = $newHandoffDateTime
